$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "N"
$ws.Range("B3").Value = "Y"
$ws.Range("B6").Value = "N"
$ws.Range("B7").Value = "Y"
$ws.Range("E4").Select()
